$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The three cells hold numeric-looking values that are stored as TEXT
# (shared strings) in the workbook. Forcing the cell to Text format before
# assigning the value keeps Excel from re-interpreting the string as a
# number, then resetting the style back to Normal avoids leaving a stray
# number-format change on the cell.

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "-0.1766"
$ws.Range("D8").Style = "Normal"

$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "-0.6006"
$ws.Range("M2").Style = "Normal"

$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "-0.5504"
$ws.Range("M3").Style = "Normal"
